# Update '想去人数' (F column) values per the diff, grouped by worksheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2808
$ws.Range("F6").Value = 2459
$ws.Range("F9").Value = 35
$ws.Range("F11").Value = 39
$ws.Range("F13").Value = 7026
$ws.Range("F14").Value = 279
$ws.Range("F15").Value = 92
$ws.Range("F16").Value = 215
$ws.Range("F18").Value = 470
$ws.Range("F19").Value = 8174
$ws.Range("F24").Value = 17
$ws.Range("F27").Value = 66
$ws.Range("F29").Value = 13
$ws.Range("F30").Value = 36
$ws.Range("F37").Value = 1152
$ws.Range("F39").Value = 646
$ws.Range("F40").Value = 3656
$ws.Range("F41").Value = 169
$ws.Range("F42").Value = 1175
$ws.Range("F43").Value = 153
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 243
$ws.Range("F15").Value = 165
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2808
$ws.Range("F5").Value = 243
$ws.Range("F6").Value = 243
$ws.Range("F8").Value = 2459
$ws.Range("F12").Value = 35
$ws.Range("F14").Value = 39
$ws.Range("F18").Value = 7026
$ws.Range("F19").Value = 279
$ws.Range("F20").Value = 92
$ws.Range("F21").Value = 215
$ws.Range("F23").Value = 470
$ws.Range("F24").Value = 8174
$ws.Range("F29").Value = 17
$ws.Range("F32").Value = 66
$ws.Range("F34").Value = 13
$ws.Range("F35").Value = 36
$ws.Range("F43").Value = 1152
$ws.Range("F44").Value = 646
$ws.Range("F45").Value = 165
$ws.Range("F46").Value = 3656
$ws.Range("F47").Value = 169
$ws.Range("F49").Value = 1175
$ws.Range("F50").Value = 153
